$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the ResellerBCN value in B8 but keep its existing formatting
$ws.Range("B8").ClearContents()

# Remove the VendorName value in F8 entirely (delete the cell content)
$ws.Range("F8").ClearContents()

# Move the active selection to B8 to match the saved selection state
$ws.Range("B8").Select()
